$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")
$notes.Range("C3").Value = "Bin barcode,`nBox barcode"

$sheetsWithHeader = @("Metadata fields for import", "sample dats", "sample cdr", "sample_open_reels")

foreach ($name in $sheetsWithHeader) {
    $ws = $wb.Worksheets.Item($name)
    if ($name -eq "Metadata fields for import") {
        $ws.Range("A2").Value = "Bin barcode"
        $ws.Range("B2").Value = "Box barcode"
    } else {
        $ws.Range("A1").Value = "Bin barcode"
        $ws.Range("B1").Value = "Box barcode"
    }
}
